# Weekly fruit/vegetable price update:
# Insert a new data row at row 24 (pushing the existing rows 24-29 down to
# 25-30) and populate it with the latest "Alcachofa" market record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24:29 down to 25:30, leaving a blank row 24 to fill in.
$ws.Rows("24:24").Insert()

$ws.Range("A24").Value2 = 1
$ws.Range("B24").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value2 = "Arica y Parinacota"
$ws.Range("D24").Value2 = 44806
$ws.Range("E24").Value2 = 15
$ws.Range("F24").Value2 = 100112013
$ws.Range("G24").Value2 = "Alcachofa"
$ws.Range("H24").Value2 = "Argentina(o)"
$ws.Range("I24").Value2 = "Primera"
$ws.Range("J24").Value2 = 250
$ws.Range("K24").Value2 = 14000
$ws.Range("L24").Value2 = 15000
$ws.Range("M24").Value2 = 14500
$ws.Range("N24").Value2 = "$/caja 40 unidades"
$ws.Range("O24").Value2 = "Provincia de Limarí"
$ws.Range("P24").Value2 = 362
$ws.Range("Q24").Value2 = 40
$ws.Range("R24").Value2 = "Hortaliza"
